# "Update Data Sources from LFX"
#
# This commit does two things to the deck:
#
#   1. Every data table that was still using the legacy custom table
#      style {AA5A75A7-0268-4E23-BE90-BAA816196550} is switched to the
#      built-in table style {48338D13-A9CA-4A3B-89B7-9073ACB4FE37}.
#
#   2. The presentation's theme colour palette is swapped: the colours
#      that used to live in the "LF Energy Theme 2023" theme become the
#      colours of the theme that PowerPoint exposes through the Design /
#      ThemeColorScheme object model (backed by ppt/theme/theme3.xml in
#      this package), while that theme's old ("Simple Light") palette is
#      retired.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style every table still on the old style GUID.
# ---------------------------------------------------------------------
$oldStyleId = "{AA5A75A7-0268-4E23-BE90-BAA816196550}"
$newStyleId = "{48338D13-A9CA-4A3B-89B7-9073ACB4FE37}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap in the "Geometric" / LF Energy colour palette as the active
#    theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# ---------------------------------------------------------------------
function ToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches the 1-12 ThemeColorScheme indices:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$newPalette = @(
    "222222",
    "FFFFFF",
    "434343",
    "999999",
    "003778",
    "0094FF",
    "5B1DE7",
    "12E2E2",
    "FF00AA",
    "ACDE1F",
    "0077CC",
    "F06292"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le 12; $k++) {
    $themeColors.Item($k).RGB = ToComRGB($newPalette[$k - 1])
}
